$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6873773333333334
$ws.Range("H2").Value = 2.062132
$ws.Range("I2").Value = 0.02660947569874856
$ws.Range("J2").Value = 0.02660947569874856
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.293362666666667
$ws.Range("N2").Value = 3.880088
$ws.Range("O2").Value = 0.02138321956175446
$ws.Range("P2").Value = 0.02138321956175446
$ws.Range("Q2").Value = 0.8890281808462223
$ws.Range("R2").Value = 8.001253627616
$ws.Range("S2").Value = 0.0005689962612895101
$ws.Range("T2").Value = 0.0005689962612895101
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6873773333333334
$ws.Range("H3").Value = 2.062132
$ws.Range("I3").Value = 0.02660947569874856
$ws.Range("J3").Value = 0.02660947569874856
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 16.09103333333333
$ws.Range("N3").Value = 48.2731
$ws.Range("O3").Value = 0.2660337333139169
$ws.Range("P3").Value = 0.2660337333139169
$ws.Range("Q3").Value = 11.06061158324444
$ws.Range("R3").Value = 99.54550424920001
$ws.Range("S3").Value = 0.007079018161664026
$ws.Range("T3").Value = 0.007079018161664026
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6873773333333334
$ws.Range("H4").Value = 2.062132
$ws.Range("I4").Value = 0.02660947569874856
$ws.Range("J4").Value = 0.02660947569874856
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 43.10054
$ws.Range("N4").Value = 129.30162
$ws.Range("O4").Value = 0.7125830471243286
$ws.Range("P4").Value = 0.7125830471243286
$ws.Range("Q4").Value = 29.62633425042667
$ws.Range("R4").Value = 266.6370082538401
$ws.Range("S4").Value = 0.01896146127579502
$ws.Range("T4").Value = 0.01896146127579502
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.913984666666667
$ws.Range("H5").Value = 17.741954
$ws.Range("I5").Value = 0.2289398029860915
$ws.Range("J5").Value = 0.2289398029860915
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.293362666666667
$ws.Range("N5").Value = 3.880088
$ws.Range("O5").Value = 0.02138321956175446
$ws.Range("P5").Value = 0.02138321956175446
$ws.Range("Q5").Value = 7.648926979105778
$ws.Range("R5").Value = 68.840342811952
$ws.Range("S5").Value = 0.004895470073676403
$ws.Range("T5").Value = 0.004895470073676403
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.913984666666667
$ws.Range("H6").Value = 17.741954
$ws.Range("I6").Value = 0.2289398029860915
$ws.Range("J6").Value = 0.2289398029860915
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 16.09103333333333
$ws.Range("N6").Value = 48.2731
$ws.Range("O6").Value = 0.2660337333139169
$ws.Range("P6").Value = 0.2660337333139169
$ws.Range("Q6").Value = 95.16212440415555
$ws.Range("R6").Value = 856.4591196374
$ws.Range("S6").Value = 0.06090571049254253
$ws.Range("T6").Value = 0.06090571049254254
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.913984666666667
$ws.Range("H7").Value = 17.741954
$ws.Range("I7").Value = 0.2289398029860915
$ws.Range("J7").Value = 0.2289398029860915
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.10054
$ws.Range("N7").Value = 129.30162
$ws.Range("O7").Value = 0.7125830471243286
$ws.Range("P7").Value = 0.7125830471243286
$ws.Range("Q7").Value = 254.8959326850534
$ws.Range("R7").Value = 2294.06339416548
$ws.Range("S7").Value = 0.1631386224198725
$ws.Range("T7").Value = 0.1631386224198725
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.230689
$ws.Range("H8").Value = 57.692067
$ws.Range("I8").Value = 0.7444507213151601
$ws.Range("J8").Value = 0.7444507213151601
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.293362666666667
$ws.Range("N8").Value = 3.880088
$ws.Range("O8").Value = 0.02138321956175446
$ws.Range("P8").Value = 0.02138321956175446
$ws.Range("Q8").Value = 24.87225520687733
$ws.Range("R8").Value = 223.850296861896
$ws.Range("S8").Value = 0.01591875322678855
$ws.Range("T8").Value = 0.01591875322678855
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.230689
$ws.Range("H9").Value = 57.692067
$ws.Range("I9").Value = 0.7444507213151601
$ws.Range("J9").Value = 0.7444507213151601
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 16.09103333333333
$ws.Range("N9").Value = 48.2731
$ws.Range("O9").Value = 0.2660337333139169
$ws.Range("P9").Value = 0.2660337333139169
$ws.Range("Q9").Value = 309.4416577219667
$ws.Range("R9").Value = 2784.9749194977
$ws.Range("S9").Value = 0.1980490046597104
$ws.Range("T9").Value = 0.1980490046597104
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.230689
$ws.Range("H10").Value = 57.692067
$ws.Range("I10").Value = 0.7444507213151601
$ws.Range("J10").Value = 0.7444507213151601
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 43.10054
$ws.Range("N10").Value = 129.30162
$ws.Range("O10").Value = 0.7125830471243286
$ws.Range("P10").Value = 0.7125830471243286
$ws.Range("Q10").Value = 828.8530804720601
$ws.Range("R10").Value = 7459.677724248541
$ws.Range("S10").Value = 0.5304829634286612
$ws.Range("T10").Value = 0.5304829634286612
